$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.133.54"
$ws.Range("E2").Value = "  -4.95%  "
$ws.Range("D3").Value = "'3.312.42"
$ws.Range("E3").Value = "  -5.20%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'569.27"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "'126.48"
$ws.Range("E6").Value = "  -5.56%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'3.311.46"
$ws.Range("E8").Value = "  -5.20%  "
$ws.Range("D9").Value = "'0.476"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").Value = "'7.16"
$ws.Range("E10").Value = "  -5.61%  "
$ws.Range("D11").Value = "'0.117"
$ws.Range("E11").Value = "  -5.75%  "
$ws.Range("D12").Value = "'0.373"
$ws.Range("E12").Value = "  -4.59%  "
$ws.Range("D13").Value = "'3.879.51"
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "'3.319.92"
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("D16").Value = "'0.0000168"
$ws.Range("E16").Value = "  -6.96%  "
$ws.Range("D17").Value = "'24.71"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "'61.217.03"
$ws.Range("E18").Value = "  -4.82%  "
$ws.Range("D19").Value = "'5.58"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").Value = "'9.06"
$ws.Range("E20").Value = "  -9.44%  "
$ws.Range("D21").Value = "'13.14"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").Value = "'351.72"
$ws.Range("E22").Value = "  -9.43%  "
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "  -4.82%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'3.444.51"
$ws.Range("E26").Value = "  -5.48%  "
$ws.Range("D27").Value = "'0.0000106"
$ws.Range("E27").Value = "  -7.60%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "'7.13"
$ws.Range("E29").Value = "  -3.39%  "
$ws.Range("D30").Value = "'1.46"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").Value = "'7.88"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("E32").Value = "  -6.80%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'0.147"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D35").Value = "'3.340.67"
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'22.28"
$ws.Range("E36").Value = "  -4.52%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'5.40"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").Value = "'6.71"
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("D39").Value = "'162.82"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("D41").Value = "'0.0753"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'40.87"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").Value = "'0.747"
$ws.Range("E44").Value = "  -7.61%  "
$ws.Range("D45").Value = "'4.23"
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").Value = "'1.11"
$ws.Range("E46").Value = "  -6.36%  "
$ws.Range("E47").Value = "  -6.69%  "
$ws.Range("D48").Value = "'22.42"
$ws.Range("E48").Value = "  -8.65%  "
$ws.Range("D49").Value = "'6.60"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("D50").Value = "'0.850"
$ws.Range("E50").Value = "  -7.51%  "
$ws.Range("D51").Value = "'2.191.70"
$ws.Range("E51").Value = "  -8.71%  "
